$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column BF holds the "Date" values for rows 2 through 31.
# The original export stamped the date text wrong (season-folder style,
# "5-17-2011-12") due to how the NBA stats were captured; correct it to a
# proper ISO date string "2012-05-17" for each data row.
#
# Excel auto-detects strings shaped like dates, so force the target cells
# to Text format first -- otherwise "2012-05-17" would be silently turned
# into a date serial number instead of staying literal text.
$rng = $ws.Range("BF2:BF31")
$rng.NumberFormat = "@"

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF$row")
    if ($cell.Value2 -eq "5-17-2011-12") {
        $cell.Value = "2012-05-17"
    }
}
